# Updates a set of odds cells in the "Jogos da Semana" workbook, as per
# the FlashScore data refresh (commit: "Atualizando o arquivo XLSX").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "J3"  = 1.04
    "K3"  = 12
    "L3"  = 1.25
    "R3"  = 1.67
    "S3"  = 2.1
    "V3"  = 9
    "W3"  = 19
    "AG3" = 12
    "AI3" = 26

    "J5" = 1.05
    "L5" = 1.3

    "J6" = 1.11
    "L6" = 1.62

    "J7" = 1.17
    "L7" = 1.73

    "G27"  = 4.8
    "H27"  = 3.6
    "N27"  = 1.87
    "O27"  = 1.75
    "Q27"  = 2.45
    "T27"  = 10
    "U27"  = 21
    "V27"  = 13
    "Z27"  = 9.25
    "AA27" = 6.2
    "AB27" = 14.5
    "AG27" = 7

    "G31"  = 1.91
    "I31"  = 3.9
    "R31"  = 1.8
    "S31"  = 1.95
    "U31"  = 9
    "AA31" = 7
    "AE31" = 12
    "AF31" = 21

    "G35"  = 6.5
    "I35"  = 1.42
    "L35"  = 1.18
    "M35"  = 4.5
    "N35"  = 1.62
    "O35"  = 2.25
    "P35"  = 1.29
    "Q35"  = 3.5
    "AA35" = 9.5
    "AE35" = 8
    "AH35" = 9.5
    "AI35" = 11

    "G42"  = 1.13
    "H42"  = 11
    "I42"  = 13
    "K42"  = 17
    "L42"  = 1.06
    "M42"  = 8
    "N42"  = 1.22
    "O42"  = 4
    "W42"  = 7.5
    "AD42" = 700
    "AE42" = 34

    "I43" = 1.42
    "N43" = 1.57
    "O43" = 2.35
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
